# Commiting exercies done in xls
# Adds a second worksheet ("Sheet2") that repeats the SD exercise from
# Sheet1 but computes the *sample* standard deviation (dividing by n-1).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update the selection on Sheet1 (it is no longer the active tab).
$ws1.Range("B2").Select()

# Insert the new worksheet right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Title
$ws2.Range("A1").Value = "Sample SD"

# Header row (bold, reuses the same labels as Sheet1)
$ws2.Range("A2").Value = "Sample"
$ws2.Range("B2").Value = "Mean"
$ws2.Range("C2").Value = "Deviation"
$ws2.Range("D2").Value = "Deviation squared"
$ws2.Range("E2").Value = "SD Deviation squared"
$ws2.Range("F2").Value = "standard deviation"
$ws2.Range("A2:F2").Font.Bold = $true

# Sample data (column A)
$ws2.Range("A3").Value = 18
$ws2.Range("A4").Value = 20
$ws2.Range("A5").Value = 23
$ws2.Range("A6").Value = 18
$ws2.Range("A7").Value = 22
$ws2.Range("A8").Value = 21
$ws2.Range("A9").Value = 17
$ws2.Range("A10").Value = 21
$ws2.Range("A11").Value = 15

# Mean (row 3 only)
$ws2.Range("B3").Formula = "=SUM(A3:A11) / 9"

# Deviation from the fixed mean constant, for every sample row
$ws2.Range("C3").Formula = "=19.44444444 - A3"
$ws2.Range("C4").Formula = "=19.44444444 - A4"
$ws2.Range("C5").Formula = "=19.44444444 - A5"
$ws2.Range("C6").Formula = "=19.44444444 - A6"
$ws2.Range("C7").Formula = "=19.44444444 - A7"
$ws2.Range("C8").Formula = "=19.44444444 - A8"
$ws2.Range("C9").Formula = "=19.44444444 - A9"
$ws2.Range("C10").Formula = "=19.44444444 - A10"
$ws2.Range("C11").Formula = "=19.44444444 - A11"

# Deviation squared
$ws2.Range("D3").Formula = "=C3^2"
$ws2.Range("D4").Formula = "=C4^2"
$ws2.Range("D5").Formula = "=C5^2"
$ws2.Range("D6").Formula = "=C6^2"
$ws2.Range("D7").Formula = "=C7^2"
$ws2.Range("D8").Formula = "=C8^2"
$ws2.Range("D9").Formula = "=C9^2"
$ws2.Range("D10").Formula = "=C10^2"
$ws2.Range("D11").Formula = "=C11^2"

# SD deviation squared (row 3 only) and standard deviation
$ws2.Range("E3").Formula = "=SUM(D3:D11)/9"
$ws2.Range("F3").Formula = "=SQRT(E3)"

# Column widths to match Sheet1's layout (closest achievable values; the
# COM width model snaps to whole-pixel increments)
$ws2.Columns("A:C").ColumnWidth = 10.666666666666666
$ws2.Columns("D").ColumnWidth = 16.333333333333336

# Selection / active cell on the new sheet
$ws2.Range("F3").Select()
